$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 'Vaga de Estágio em Arquitetura e Urbanismo'
$ws.Range("C2").Value = 'SUPER ESTAGIOS'
$ws.Range("D2").Value = 'Auxiliar nas seguintes atividades Acompanhamento de obras Suporte na Elaboração de projetos Experiência de campo Levantamento arquitetônico Projetos interiores Necessário Excel intermediário AutoCad e Reddit A partir do 7 período em um dos cursos Morar próximo a Jacarepaguá'
$ws.Range("E2").Value = 200000
$ws.Range("F2").Value = 'Rio de Janeiro - RJ (1)'
$ws.Range("G2").Value = 'https://www.catho.com.br/vagas/estagio-em-arquitetura-e-urbanismo/30692595/'

$ws.Range("D3").Value = 'Realizar lançamentos de pagamentos no sistema auxiliar na organização da Unidade agenda e manutenção de agendas de reuniões acompanhar prestadores de serviço fazer as cotações e compras de materiais de insumo da Unidade
 Ter boa comunicação Ser proativo'
$ws.Range("E3").Value = 200000
$ws.Range("F3").Value = 'São Paulo - SP (1)'
$ws.Range("G3").Value = 'https://www.catho.com.br/vagas/estagio/30731456/'

$ws.Range("B4").Value = 'Vaga de Estágio Fonoaudiologia'
$ws.Range("C4").Value = 'CLÍNICA MENTHALHELP'
$ws.Range("D4").Value = 'O Fonoaudiólogo é responsável por avaliar pacientes com distúrbios de comunicação identificando problemas de fala linguagem voz audição e motricidade oral Desenvolvimento de planos de tratamento Com base na avaliação dos pacientes o Fonoaudiólogo desenvolve planos de tratamento personalizados com o objetivo de ajudar os pacientes a superar seus problemas de comunicação Realização de terapias O Fonoaudiólogo realiza terapias para ajudar os pacientes a desenvolver habilidades de comunicação corrigir problemas de fala linguagem voz audição e motricidade oral Acompanhamento e monitoramento O Fonoaudiólogo acompanha e monitora o progresso dos pacientes durante o tratamento realizando ajustes nos planos de tratamento quando necessário Orientação aos pacientes e familiares O Fonoaudiólogo orienta os pacientes e seus familiares sobre como melhorar a comunicação e prevenir futuros distúrbios Trabalho em equipe O Fonoaudiólogo trabalha em equipe com outros profissionais da saúde como médicos psicólogos e terapeutas ocupacionais para fornecer tratamento integrado aos pacientes Pesquisa e educação Alguns Fonoaudiólogos realizam pesquisas e participam de atividades de educação continuada para manterse atualizados sobre os avanços em sua área e melhorar suas habilidades profissionais O objetivo geral do trabalho do Fonoaudiólogo é melhorar a qualidade de vida dos pacientes ajudandoos a superar seus problemas de comunicação e melhorar sua autoestima e autoconfiança O trabalho do Fonoaudiólogo é importante em muitos contextos incluindo escolas hospitais clínicas e empresas onde pode ajudar a melhorar a comunicação e a produtividade das pessoas
Ter fácil acesso a região de Guarulhos
Deverá estar cursando Fonoaudiologia Estar cursando mínimo 4 semestre Estar cursando fonoaudiologia mínimo 4 semestre'
$ws.Range("E4").Value = 200100300000
$ws.Range("F4").Value = 'Guarulhos - SP (2)'
$ws.Range("G4").Value = 'https://www.catho.com.br/vagas/estagio-fonoaudiologia/30078566/'

$ws.Range("B5").Value = 'Vaga de Estágio - Financeiro'
$ws.Range("C5").Value = 'ELIS ENERGIA'
$ws.Range("D5").Value = '1 Apoio na análise de dados financeiros e contábeis 
2 Apoio no desenvolvimento de apresentações e memorandos
3 Auxílio nas pesquisas de mercado e players 
4 Auxílio na elaboração de modelos financeiros para avaliação de projetos e colaborações em teses de investimento e oportunidades de MA
5 Apoio na elaboração de relatórios financeiros e contábeis Boa comunicação interpessoal e trabalho em equipe Crítico Senso de Dono Atitude proativa e capacidade de resolver problemas de forma eficaz Interesse em aprender sobre o contexto financeiro em empresa de energia solar
Conhecimento Intermediário do Pacote Office Excel e PowerPoint Avançado 
Inglês Intermediário será um diferencial'
$ws.Range("E5").Value = 200000
$ws.Range("F5").Value = 'São Paulo - SP (1)'
$ws.Range("G5").Value = 'https://www.catho.com.br/vagas/estagio-financeiro/30730987/'

$ws.Range("B6").Value = 'Vaga de Estágio em Direito'
$ws.Range("C6").Value = 'ADVOCACIA FELIZARDO BARROSO & ASSOCIADOS'
$ws.Range("D6").Value = 'Acompanhamento de Processos Elaboração de Petições Distribuições Extração de Guias Cursando direito entre o 4 e 8 períodos'
$ws.Range("E6").Value = 200000
$ws.Range("F6").Value = 'Rio de Janeiro - RJ (1)'
$ws.Range("G6").Value = 'https://www.catho.com.br/vagas/estagio-em-direito/30733545/'

$ws.Range("B7").Value = 'Vaga de Estágio de Marketing'
$ws.Range("C7").Value = 'EMPRESA CONFIDENCIAL'
$ws.Range("D7").Value = ' Atuar no apoio às atividades de marketing 
 Atuar no apoio aos eventos 
 Atuar com demandas na organização de Podcast'
$ws.Range("E7").Value = 200000
$ws.Range("F7").Value = 'Cuiaba - MT (1)'
$ws.Range("G7").Value = 'https://www.catho.com.br/vagas/estagio-de-marketing/30731568/'

$ws.Range("B8").Value = 'Vaga de Estágio de Enfermagem'
$ws.Range("C8").Value = 'COMERCIAL DAHANA LIMITADA'
$ws.Range("D8").Value = 'Realizar avaliações de saúde ocupacional
Realizar campanhas de saúde
Realizar atendimento ambulatorial
Realizar arquivo de documentos
Documentar registros de saúde Deverá estar cursando Graduação em andamento no curso de enfermagem
Disponibilidade para atuar em Contagem
Disponibilidade para estagiar de 0800 às 1500'
$ws.Range("E8").Value = 200000
$ws.Range("F8").Value = 'Contagem - MG (1)'
$ws.Range("G8").Value = 'https://www.catho.com.br/vagas/estagio-de-enfermagem/30731495/'

$ws.Range("B9").Value = 'Vaga de Estágio em Administração'
$ws.Range("C9").Value = 'ADVOCACIA FELIZARDO BARROSO & ASSOCIADOS'
$ws.Range("D9").Value = 'Atuar na área administrativa'
$ws.Range("E9").Value = 200000
$ws.Range("F9").Value = 'Rio de Janeiro - RJ (1)'
$ws.Range("G9").Value = 'https://www.catho.com.br/vagas/estagio-em-administracao/30733464/'

$ws.Range("B10").Value = 'Vaga de Estágio na área Administrativa'
$ws.Range("C10").Value = 'EMPRESA CONFIDENCIAL'
$ws.Range("D10").Value = ' Auxílio na distribuição de demandas de estagiários  Atualização de Planilhas de Controle de Qualidade  Auxílio no Tratamento de Dados em Planilhas do Google Sheets  Acompanhamento de Gestão de Estagiários sanando dúvidas de escala envio de atestados e outros  Revisão de textos
  Boa Comunicação oral e escrita  Boa organização  Conhecimento linguístico gramatical e de norma culta  Se interessar por escrever e revisar texto'
$ws.Range("E10").Value = 200000
$ws.Range("F10").Value = 'Campinas - SP (1)'
$ws.Range("G10").Value = 'https://www.catho.com.br/vagas/estagio-na-area-administrativa/30731099/'

$ws.Range("B11").Value = 'Vaga de Estagio em contabilidade'
$ws.Range("C11").Value = 'HOMMAGE CONTABILIDADE'
$ws.Range("D11").Value = 'Escrituração de notas fiscais de entrada e saída mercadorias e serviços geração de boletos do escritório e entregar em contato com os clientes para realizar cobranças fechando os recebíveis
 Ser pró ativo ser responsável pelo trabalho cumprindo metas e prazos ter boa comunicação interna e externa'
$ws.Range("E11").Value = 200000
$ws.Range("F11").Value = 'São Paulo - SP (1)'
$ws.Range("G11").Value = 'https://www.catho.com.br/vagas/estagio-em-contabilidade/30732901/'

$ws.Range("B12").Value = 'Vaga de Estágio em Projetos'
$ws.Range("C12").Value = 'FERRARI SOLUÇÕES EM ENGENHARIA'
$ws.Range("D12").Value = 'Carga horária 4 horas  Período a combinar
Auxiliar na elaboração de projetos das edificações com plantas cortes e fachadas medições e conferência de layout in loco
Auxiliar na implantação do sistema de segurança contra incêndio em projetos memoriais de cálculos e documentos
Auxiliar nas vistorias técnicas para adequações de imóveis no âmbito da segurança contra incêndio
Auxiliar elaboração de orçamentos atendimento telefônico rotinas administrativas gerenciamento de compras e logística de materiais
  Estagiário de curso arquiteturaengenhariatecnólogo partir do terceiro semestre Boa escrita e interpretação de textos e projetos
Conhecimentos de projeto avançado em Auto CAD
 Conhecimentos no pacote Office Word Outlook Excel desejável conhecimento em maquete 3D
 Fácil acesso estamos localizados próximo a linha amarela estação Vila Sônia'
$ws.Range("E12").Value = 200000
$ws.Range("F12").Value = 'São Paulo - SP (1)'
$ws.Range("G12").Value = 'https://www.catho.com.br/vagas/estagio-em-projetos/30733327/'

$ws.Range("B13").Value = 'Vaga de VAGA DE ESTÁGIO EM DIREITO'
$ws.Range("C13").Value = 'IDEALIZA CIDADES'
$ws.Range("D13").Value = ' Elaboração de contratos
 Pesquisa de jurisprudência
 Análise de legislação 
 Auxiliar na elaboração de pareceres jurídicos
 Auxiliar no processo de registro de loteamentos e condomínios
 Elaboração de documentos para registro de Associação de Moradores
 Cursando 3º ano de Direito
Disponibilidade de estágio presencial
Interesse na área do Direito Imobiliário'
$ws.Range("E13").Value = 200000
$ws.Range("F13").Value = 'São Paulo - SP (1)'
$ws.Range("G13").Value = 'https://www.catho.com.br/vagas/vaga-de-estagio-em-direito/30732830/'

$ws.Range("B14").Value = 'Vaga de VAGA DE ESTÁGIO EM DIREITO'
$ws.Range("C14").Value = 'IDEALIZA CIDADES'
$ws.Range("D14").Value = ' Elaboração de contratos
 Pesquisa de jurisprudência
 Análise de legislação 
 Auxiliar na elaboração de pareceres jurídicos
 Auxiliar no processo de registro de loteamentos e condomínios
 Elaboração de documentos para registro de Associação de Moradores
 Ï Cursando 3º ano de Direito
Ï Disponibilidade de estágio presencial
Ï Interesse na área do Direito Imobiliário'
$ws.Range("E14").Value = 200000
$ws.Range("G14").Value = 'https://www.catho.com.br/vagas/vaga-de-estagio-em-direito/30731462/'

$ws.Range("B15").Value = 'Vaga de Estágio Administrativo'
$ws.Range("C15").Value = 'ELIS BRASIL'
$ws.Range("D15").Value = 'Atuará em atividades relacionadas a área administrativa e atendimento de clientes internos e externos Residir em Jundiaí e região'
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 'Jundiai - SP (1)'
$ws.Range("G15").Value = 'https://www.catho.com.br/vagas/estagio-administrativo/30730722/'

$ws.Range("B16").Value = 'Vaga de Estágio em Economia'
$ws.Range("C16").Value = 'EMPRESA CONFIDENCIAL'
$ws.Range("D16").Value = 'conciliação e lançamento das operações dos fundos de investimento execução das rotinas operacionais dos sistemas contratados suporte na elaboração de relatórios diários e mensais de rentabilidade e risco das carteiras gerenciamento de caixa Excel'
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 'Diadema - SP (1)'
$ws.Range("G16").Value = 'https://www.catho.com.br/vagas/estagio-em-economia/30732579/'

$ws.Range("B17").Value = 'Vaga de Estágio em Vendas'
$ws.Range("C17").Value = 'EVOCONT CONTABILIDADE DIGITAL'
$ws.Range("D17").Value = 'Atividades
Atuação junto a equipe comercial no atendimento de cliente em loja
Auxiliar no controle de estoque físico e fluxo de caixa
Auxiliar na geração de orçamentos e cadastro de clientes
Auxiliar em pedidos para abastecimento de estoque
Auxiliar na divulgação de campanhas nas redes sociais
Auxiliar na organização geral do mostruário
Habilidades
Ótima comunicação proatividade e organização
Requisitos
Cursando Gestão Comercial Processos Gerenciais Administração Marketing e áreas correlatas
Conhecimento intermediário do Pacote Office
Jornada de Trabalho
Segunda a sextafeira 0900h às 1500h '
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 'Mogi-guaçu - SP (1)'
$ws.Range("G17").Value = 'https://www.catho.com.br/vagas/estagio-em-vendas/30732770/'

$ws.Range("B18").Value = 'Estágio administrativo - Contas a receber'
$ws.Range("C18").Value = 'Estágio administrativo - Contas a receber'
$ws.Range("D18").Value = 'DESCRIÇÃO DA VAGA
Você quer fazer parte da maior especialista em tecnologia da América Latina Com um portfólio que reúne mais de 10000 fabricantes e mais de 1000000 de títulos em software estamos presentes no Brasil México e Colômbia com um alcance que se estende por toda a América
Se você é apaixonado por tecnologia e busca constante aprendizado queremos você em nosso time
Venha trilhar uma jornada de crescimento e desenvolvimento profissional que fará a diferença na sua carreira
Inscrevase agora e ajude a vamos juntos construir o futuro da tecnologia'
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 'São Paulo - SP'
$ws.Range("G18").Value = 'https://software.gupy.io/job/eyJqb2JJZCI6ODAxNTEwNiwic291cmNlIjoiZ3VweV9wb3J0YWwifQ==?jobBoardSource=gupy_portal'

$ws.Range("B19").Value = 'Estágio em Engenharia Civil'
$ws.Range("C19").Value = 'Estágio em Engenharia Civil'
$ws.Range("D19").Value = 'DESCRIÇÃO DA VAGA
Antes de falarmos sobre a vaga que tal entender um pouco sobre nós
Somos uma construtora pura focada em empreendimentos de médio e alto padrão na região metropolitana de São Paulo Temos como propósito desenvolver espaços onde as pessoas sejam acolhidas com dignidade prazer e se sintam melhores onde moram trabalham se hospedam consomem e investem
São mais de 40 obras concluídas Em nosso portifólio temos obras residenciais comerciais hotéis e saúde
Nossos pilares de negócios são a excelência a sustentabilidade e a inovação
Valorizamos nossos talentos por isso mais de 40 do nosso time é formado dentro de casa Estamos crescendo muito e claro queremos que nosso time cresça com a gente
Ajudar a construir o futuro de São Paulo não é somente criar projetos é também gerar valor pra sociedade através da solidariedade de projetos sustentáveis parcerias duradouras e respeito ao meio ambiente Buscamos quem partilha destes princípios e acredita ser possível transformar o futuro de São Paulo
Criamos um ambiente familiar e profissional de integração que proporcione a coletividade de reconhecimento mútuo estimulamos que seja colaborativo e positivo Aqui todas as pessoas são bem vindas e respeitadas em suas diferentes origens crenças experiências raças deficiências orientações sexuais e gerações
VEMPRAROCONTEC'
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 'São Paulo - SP'
$ws.Range("G19").Value = 'https://rocontec.gupy.io/job/eyJqb2JJZCI6ODAxMTA2Mywic291cmNlIjoiZ3VweV9wb3J0YWwifQ==?jobBoardSource=gupy_portal'

$ws.Range("B20").Value = 'Estagio em Engenharia Civil - Obra'
$ws.Range("C20").Value = 'Estagio em Engenharia Civil - Obra'
$ws.Range("D20").Value = 'DESCRIÇÃO DA VAGA
Antes de falarmos sobre a vaga que tal entender um pouco sobre nós
Somos uma construtora pura focada em empreendimentos de médio e alto padrão na região metropolitana de São Paulo Temos como propósito desenvolver espaços onde as pessoas sejam acolhidas com dignidade prazer e se sintam melhores onde moram trabalham se hospedam consomem e investem
São mais de 40 obras concluídas Em nosso portifólio temos obras residenciais comerciais hotéis e saúde
Nossos pilares de negócios são a excelência a sustentabilidade e a inovação
Valorizamos nossos talentos por isso mais de 40 do nosso time é formado dentro de casa Estamos crescendo muito e claro queremos que nosso time cresça com a gente
Ajudar a construir o futuro de São Paulo não é somente criar projetos é também gerar valor pra sociedade através da solidariedade de projetos sustentáveis parcerias duradouras e respeito ao meio ambiente Buscamos quem partilha destes princípios e acredita ser possível transformar o futuro de São Paulo
Criamos um ambiente familiar e profissional de integração que proporcione a coletividade de reconhecimento mútuo estimulamos que seja colaborativo e positivo Aqui todas as pessoas são bem vindas e respeitadas em suas diferentes origens crenças experiências raças deficiências orientações sexuais e gerações
VEMPRAROCONTEC'
$ws.Range("E20").Value = 0
$ws.Range("G20").Value = 'https://rocontec.gupy.io/job/eyJqb2JJZCI6Nzk1NTE1MSwic291cmNlIjoiZ3VweV9wb3J0YWwifQ==?jobBoardSource=gupy_portal'

$ws.Range("B21").Value = 'Programa de Estágio - Futuros Líderes Maxpar'
$ws.Range("C21").Value = 'Programa de Estágio - Futuros Líderes Maxpar'
$ws.Range("D21").Value = 'DESCRIÇÃO DA VAGA
Olá D
 Já pensou em fazer parte de um Grupo Brasileiro com mais de 6000 colaboradores e ter a oportunidade conhecer diferentes áreas além de ser capacitado a e desafiado a e se tornar um futuro líder de equipe
 Estamos em busca de talentos que se identifiquem com uma carreira de liderança que tenham sede de aprendizado e queiram fazer a diferença O programa Futuros Líderes tem duração de 6 meses a 1 ano e o objetivo vai muito além da efetivação de um contrato de estágio queremos formar líderes para crescer com a gente
 Você terá a oportunidade de atuar e conhecer diversos setores relacionados ao Negócio Maxpar e receber suporte e mentoria de gestores que são referência fornecendo um processo de desenvolvimento contínuo e enriquecedor
 Se você é apaixonadoa por inovação pensa fora da caixa trazendo soluções busca um ambiente de ebulição de ideias gosta de trabalhar em equipe deseja desenvolvimento contínuo e plano de carreira seu lugar é aqui Inscrevase hoje 
 FuturosLíderesMaxpar
OrgulhoDeSerGrupoAutoglass
 '
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 'Vila Velha - ES'
$ws.Range("G21").Value = 'https://autoglassestagio.gupy.io/job/eyJqb2JJZCI6ODAxNDUyMiwic291cmNlIjoiZ3VweV9wb3J0YWwifQ==?jobBoardSource=gupy_portal'

$ws.Range("B22").Value = 'Pessoa Estagiária de Atendimento'
$ws.Range("C22").Value = 'Pessoa Estagiária de Atendimento'
$ws.Range("D22").Value = 'DESCRIÇÃO DA VAGA
A SiMCo foi criada com o foco em oferecer atendimento médico e odontológico de qualidade aos 150 milhões de brasileiros que não possuem planos de saúde privados SiM significa Serviço de Inclusão à Medicina e a cor verde sinaliza para as pessoas sinal verde para um serviço de saúde acessível de alta qualidade
Vem com a gente construir a maior e melhor plataforma de acesso à saúde da América Latina '
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 'Fortaleza'
$ws.Range("G22").Value = 'https://simco.gupy.io/job/eyJqb2JJZCI6Nzk4OTE4Niwic291cmNlIjoiZ3VweV9wb3J0YWwifQ==?jobBoardSource=gupy_portal'

$ws.Range("B23").Value = 'Estágio - Trading de Energia'
$ws.Range("C23").Value = 'Estágio - Trading de Energia'
$ws.Range("D23").Value = 'DESCRIÇÃO DA VAGA
Nossa cliente é uma empresa pioneira no Brasil autorizada a comercializar energia elétrica com consumidores finais e geradores no ambiente de contratação livre
Área Comercialização de Energia'
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 'Curitiba - PR'
$ws.Range("G23").Value = 'https://vagascetefe.gupy.io/job/eyJqb2JJZCI6ODAwNDk3Miwic291cmNlIjoiZ3VweV9wb3J0YWwifQ==?jobBoardSource=gupy_portal'

$ws.Range("B24").Value = 'Estágio - Atendimento'
$ws.Range("C24").Value = 'Estágio - Atendimento'
$ws.Range("D24").Value = 'DESCRIÇÃO DA VAGA
Buscamos uma pessoa para atuar como estagiário de atendimento na MRM Brasil uma agência fullservice com maior foco em digital do McCann Worldgroup unidade do Grupo Interpublic IPG Uma agência que acredita no poder das intersecções entre estratégia criatividade e tecnologia suportado por dados para criar experiências totalmente humanas Atuamos a partir do modelo híbrido de agência e consultoria de transformação digital ajudando a fortalecer o relacionamento entre negócios e pessoas ajudando marcas a encontrarem seu propósito aprofundar suas relações e criar experiências que gerem crescimento
 A MRM Brasil valoriza a criatividade e inovação algo que é potencializado em um ambiente diverso e inclusivo por isso damos prioridade na contratação de grupos minoritários como pessoas com mais de 45 anos pretas pardas indígenas pessoas da comunidade LGBTQIAP e PCDs
 Venha fazer parte de um time Better Together onde sua individualidade é respeitada e valorizada em todos os sentidos'
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 'São Paulo - SP'
$ws.Range("G24").Value = 'https://mrm.gupy.io/job/eyJqb2JJZCI6ODAxMzIyNSwic291cmNlIjoiZ3VweV9wb3J0YWwifQ==?jobBoardSource=gupy_portal'

$ws.Range("B25").Value = 'Estágio em Segurança do Trabalho'
$ws.Range("C25").Value = 'Estágio em Segurança do Trabalho'
$ws.Range("D25").Value = 'DESCRIÇÃO DA VAGA
A Krones do Brasil Ltda está em busca de um profissional motivado que almeje se desenvolver profissionalmente e queira crescer dentro da empresa para fazer parte da nossa equipe de Segurança do Trabalho O profissional executará atividades para eliminação de riscos e prevenção de acidentes na organização visando proteger e resguardar a integridade dos empregados próprios e terceiros bem como apoiará os de Meio Ambiente e Qualidade
SolutionsBeyondTomorrow BePartOfKrones'
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 'Diadema - SP'
$ws.Range("G25").Value = 'https://krones.gupy.io/job/eyJqb2JJZCI6ODAxNDEzMSwic291cmNlIjoiZ3VweV9wb3J0YWwifQ==?jobBoardSource=gupy_portal'

$ws.Range("B26").Value = 'Estágio de Social Media'
$ws.Range("C26").Value = 'Estágio de Social Media'
$ws.Range("D26").Value = 'DESCRIÇÃO DA VAGA
Se você se identifica com uma gestão horizontal com menos burocracia e mais ação e valoriza a autonomia para trazer e concretizar ideias evoluindo e crescendo no processo venha fazer parte do nosso time de uniques
Estamos em busca de uma pessoa para atuar como estagiária de Social Media na Briefing
A Briefing é nossa agência de DJs e bandas com a missão de identificar as necessidades de diferentes tipos de clientes e oferecer o artista mais adequado
A vaga é presencial de segunda a quinta com home office às sextas no período da manhã ou tarde
A empresa está localizada no Jardim Paulista próxima ao Parque Ibirapuera e à Avenida Juscelino Kubitschek'
$ws.Range("E26").Value = 0
$ws.Range("G26").Value = 'https://umauma.gupy.io/job/eyJqb2JJZCI6ODAxMzcxMywic291cmNlIjoiZ3VweV9wb3J0YWwifQ==?jobBoardSource=gupy_portal'

